$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 18-23: tutoring session dates (col B) and durations (col C)
$ws.Range("B18").Value2 = 43987
$ws.Range("C18").Value = "3h"
$ws.Range("B19").Value2 = 43994
$ws.Range("C19").Value = "3h"
$ws.Range("B20").Value2 = 44014
$ws.Range("C20").Value = "2h"
$ws.Range("B21").Value2 = 44031
$ws.Range("C21").Value = "2h"
$ws.Range("B22").Value2 = 44046
$ws.Range("C22").Value = "2h"
$ws.Range("B23").Value2 = 44053

# Apply a short-date number format to B18, then copy just the formatting
# down the rest of the date column so every cell shares a single style.
$ws.Range("B18").NumberFormat = "mm-dd-yy"
$ws.Range("B18").Copy()
$ws.Range("B19:B23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 17: "Tutorias" section header - bold white text centered on a black fill
$hdr = $ws.Range("B17")
$hdr.Value = "Tutorias"
$hdr.Font.Bold = $true
$hdr.Font.ThemeColor = 2
$hdr.Interior.ThemeColor = 1
$hdr.HorizontalAlignment = -4108

# Row 15: note about the next tutoring session
$ws.Range("B15").Value = "Siguiente turoría día 12 16:00h"

# Scroll the view and move the selection to match the saved window state
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$ws.Range("D20").Select()
